{"js": "// Replace the date line and every three-digit-by-one-digit multiplication\n// answer in the worksheet table. Every \"before\" string below occurs exactly\n// once in the document, so a scoped, case-sensitive search-and-replace for\n// each pair (applied in order) unambiguously reproduces the target edit,\n// including the pair at index 17/19 that momentarily reuses the same text\n// (\"804\u00d73=2412\") for two different cells.\nconst replacements = [\n  [\"2024-06-16 Sunday\", \"2024-06-17 Monday\"],\n  [\"471\u00d78=3768\", \"628\u00d78=5024\"],\n  [\"388\u00d73=1164\", \"814\u00d75=4070\"],\n  [\"638\u00d73=1914\", \"369\u00d76=2214\"],\n  [\"429\u00d75=2145\", \"574\u00d77=4018\"],\n  [\"324\u00d77=2268\", \"680\u00d72=1360\"],\n  [\"474\u00d78=3792\", \"420\u00d79=3780\"],\n  [\"671\u00d76=4026\", \"762\u00d73=2286\"],\n  [\"262\u00d78=2096\", \"518\u00d73=1554\"],\n  [\"324\u00d78=2592\", \"434\u00d78=3472\"],\n  [\"117\u00d75=585\", \"924\u00d76=5544\"],\n  [\"278\u00d78=2224\", \"281\u00d79=2529\"],\n  [\"580\u00d73=1740\", \"293\u00d78=2344\"],\n  [\"128\u00d72=256\", \"167\u00d73=501\"],\n  [\"788\u00d75=3940\", \"700\u00d72=1400\"],\n  [\"784\u00d78=6272\", \"381\u00d78=3048\"],\n  [\"526\u00d74=2104\", \"998\u00d72=1996\"],\n  [\"804\u00d73=2412\", \"815\u00d72=1630\"],\n  [\"275\u00d75=1375\", \"310\u00d73=930\"],\n  [\"486\u00d72=972\", \"804\u00d73=2412\"],\n  [\"430\u00d79=3870\", \"115\u00d74=460\"],\n  [\"843\u00d72=1686\", \"778\u00d73=2334\"],\n  [\"171\u00d72=342\", \"547\u00d72=1094\"],\n  [\"843\u00d73=2529\", \"914\u00d76=5484\"],\n  [\"418\u00d77=2926\", \"525\u00d78=4200\"],\n  [\"626\u00d76=3756\", \"628\u00d79=5652\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for ${JSON.stringify(oldText)}, found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each \"before\" string occurs exactly once in the document at the time it is\n# looked up (verified against the source diff), including the pair that\n# temporarily reuses \"804x3=2412\" for two different cells -- so replacing in\n# this exact order, one occurrence at a time, reproduces the target edit.\n$pairs = @(\n  ,@(\"2024-06-16 Sunday\", \"2024-06-17 Monday\")\n  ,@(\"471\u00d78=3768\", \"628\u00d78=5024\")\n  ,@(\"388\u00d73=1164\", \"814\u00d75=4070\")\n  ,@(\"638\u00d73=1914\", \"369\u00d76=2214\")\n  ,@(\"429\u00d75=2145\", \"574\u00d77=4018\")\n  ,@(\"324\u00d77=2268\", \"680\u00d72=1360\")\n  ,@(\"474\u00d78=3792\", \"420\u00d79=3780\")\n  ,@(\"671\u00d76=4026\", \"762\u00d73=2286\")\n  ,@(\"262\u00d78=2096\", \"518\u00d73=1554\")\n  ,@(\"324\u00d78=2592\", \"434\u00d78=3472\")\n  ,@(\"117\u00d75=585\", \"924\u00d76=5544\")\n  ,@(\"278\u00d78=2224\", \"281\u00d79=2529\")\n  ,@(\"580\u00d73=1740\", \"293\u00d78=2344\")\n  ,@(\"128\u00d72=256\", \"167\u00d73=501\")\n  ,@(\"788\u00d75=3940\", \"700\u00d72=1400\")\n  ,@(\"784\u00d78=6272\", \"381\u00d78=3048\")\n  ,@(\"526\u00d74=2104\", \"998\u00d72=1996\")\n  ,@(\"804\u00d73=2412\", \"815\u00d72=1630\")\n  ,@(\"275\u00d75=1375\", \"310\u00d73=930\")\n  ,@(\"486\u00d72=972\", \"804\u00d73=2412\")\n  ,@(\"430\u00d79=3870\", \"115\u00d74=460\")\n  ,@(\"843\u00d72=1686\", \"778\u00d73=2334\")\n  ,@(\"171\u00d72=342\", \"547\u00d72=1094\")\n  ,@(\"843\u00d73=2529\", \"914\u00d76=5484\")\n  ,@(\"418\u00d77=2926\", \"525\u00d78=4200\")\n  ,@(\"626\u00d76=3756\", \"628\u00d79=5652\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n  if (-not $found) {\n    throw \"Could not find text: $oldText\"\n  }\n}\n\nWrite-Output \"done\""}
